$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.427.09"
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = "'1.869.16"
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'247.40"
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = "'0.4730"
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").Value = "'0.2912"
$ws.Range("E8").Value = '  +1.68%  '
$ws.Range("D9").Value = "'0.06478"
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = "'22.03"
$ws.Range("E10").Value = '  +5.95%  '
$ws.Range("D11").Value = "'0.07720"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = "'97.55"
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").Value = "'0.7404"
$ws.Range("E13").Value = '  +4.89%  '
$ws.Range("D14").Value = "'1.871.55"
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = "'5.142"
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("D16").Value = "'273.09"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = "'30.409.98"
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").Value = "'13.39"
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("D20").Value = "'0.000007499"
$ws.Range("E20").Value = '  -0.57%  '
$ws.Range("D21").Value = "'2.115.79"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = "'5.234"
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").Value = "'6.166"
$ws.Range("D25").Value = "'9.268"
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = "'163.39"
$ws.Range("E26").Value = '  -1.02%  '
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = "'1.917"
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").Value = "'0.09994"
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("D30").Value = "'1.364"
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("D31").Value = "'1.507"
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = "'4.115"
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("D34").Value = "'0.04807"
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("D35").Value = "'1.115"
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = "'0.6932"
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("D38").Value = "'0.01848"
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("E39").Value = '  +0.40%  '
$ws.Range("D40").Value = "'6.266"
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("E41").Value = '  +3.88%  '
$ws.Range("D42").Value = "'1.968"
$ws.Range("E42").Value = '  +3.98%  '
$ws.Range("D43").Value = "'0.4182"
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = "'0.8355"
$ws.Range("E45").Value = '  -0.34%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = "'9.288"
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = "'35.46"
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").Value = "'6.962"
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = "'920.40"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").Value = "'0.05637"
$ws.Range("E51").Value = '  +1.34%  '
